# Updates the "Estado de Cuenta" workbook:
#  - removes the EUDIN HERNANDEZ MORALES worker block (2 rows)
#  - re-sorts the remaining worker periods in ascending order
#  - refreshes the "Valor Mora" period amounts to the new figures
#  - updates the summary totals (Cant. Trabajadores / Cant. Periodos / Valor Mora)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- remove the two rows that held EUDIN HERNANDEZ MORALES (CC 12602950) ---
$ws.Rows.Item(16).EntireRow.Delete()
$ws.Rows.Item(16).EntireRow.Delete()

# --- rewrite the worker detail rows (now at 16-24) in period-ascending order ---
$rows = @(
    @{ Row = 16; Tipo = "PE"; Doc = "800176329091986"; Nombre = "ABRAHAM JOSE GIL HERNANDEZ"; Periodo = "2011"; Valor = 7022;  Salario = 877803 },
    @{ Row = 17; Tipo = "PE"; Doc = "800176329091986"; Nombre = "ABRAHAM JOSE GIL HERNANDEZ"; Periodo = "2012"; Valor = 35112; Salario = 877803 },
    @{ Row = 18; Tipo = "PE"; Doc = "800176329091986"; Nombre = "ABRAHAM JOSE GIL HERNANDEZ"; Periodo = "2101"; Valor = 35112; Salario = 877803 },
    @{ Row = 19; Tipo = "PE"; Doc = "800176329091986"; Nombre = "ABRAHAM JOSE GIL HERNANDEZ"; Periodo = "2102"; Valor = 35112; Salario = 877803 },
    @{ Row = 20; Tipo = "CC"; Doc = "1082471193";      Nombre = "JOSE LUIS CARRASCAL MACHADO"; Periodo = "2103"; Valor = 36341; Salario = 908526 },
    @{ Row = 21; Tipo = "CC"; Doc = "1082471193";      Nombre = "JOSE LUIS CARRASCAL MACHADO"; Periodo = "2104"; Valor = 36341; Salario = 908526 },
    @{ Row = 22; Tipo = "CC"; Doc = "1082471193";      Nombre = "JOSE LUIS CARRASCAL MACHADO"; Periodo = "2105"; Valor = 36341; Salario = 908526 },
    @{ Row = 23; Tipo = "CC"; Doc = "1082471193";      Nombre = "JOSE LUIS CARRASCAL MACHADO"; Periodo = "2106"; Valor = 36341; Salario = 908526 },
    @{ Row = 24; Tipo = "CC"; Doc = "1082471193";      Nombre = "JOSE LUIS CARRASCAL MACHADO"; Periodo = "2107"; Valor = 19382; Salario = 908526 }
)

foreach ($r in $rows) {
    $ws.Range("B" + $r.Row).Value = $r.Tipo
    $ws.Range("C" + $r.Row).Value = $r.Doc
    $ws.Range("D" + $r.Row).Value = $r.Nombre
    $ws.Range("E" + $r.Row).Value = $r.Periodo
    $ws.Range("F" + $r.Row).Value = $r.Valor
    $ws.Range("G" + $r.Row).Value = $r.Salario
}

# --- update the summary block ---
$ws.Range("E11").Value = 277104
$ws.Range("C13").Value = 2
$ws.Range("F13").Value = 9
